$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MuSCs -> FAPs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5854969999999999
$ws.Range("H2").Value = 1.170994
$ws.Range("O2").Value = 0.4086672402490986
$ws.Range("P2").Value = 0.5089958879585649
$ws.Range("Q2").Value = 0.023360744803
$ws.Range("R2").Value = 0.140164468818
$ws.Range("S2").Value = 0.4086672402490986
$ws.Range("T2").Value = 0.5089958879585649

# Row 3 (MuSCs -> MuSCs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5854969999999999
$ws.Range("H3").Value = 1.170994
$ws.Range("M3").Value = 0.057733
$ws.Range("N3").Value = 0.115466
$ws.Range("O3").Value = 0.5913327597509014
$ws.Range("P3").Value = 0.4910041120414351
$ws.Range("Q3").Value = 0.033802498301
$ws.Range("R3").Value = 0.135209993204
$ws.Range("S3").Value = 0.5913327597509014
$ws.Range("T3").Value = 0.4910041120414351
